$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$cell = $t.Cell(3, 2)
$r = $d.Range($cell.Range.Start, $cell.Range.End)
$r.Find.Execute("0.7 (0.7, 0.8)", $true, $false, $false, $false, $false, $true, 0, $false, "0.741 (0.705, 0.776)", 2)
$cell = $t.Cell(3, 3)
$r = $d.Range($cell.Range.Start, $cell.Range.End)
$r.Find.Execute("4.8 (4.7, 4.9)", $true, $false, $false, $false, $false, $true, 0, $false, "4.817 (4.727, 4.907)", 2)
$cell = $t.Cell(3, 4)
$r = $d.Range($cell.Range.Start, $cell.Range.End)
$r.Find.Execute("13.0 (10.2, 16.0)*", $true, $false, $false, $false, $false, $true, 0, $false, "13.046 (10.204, 15.960)*", 2)
$cell = $t.Cell(3, 6)
$r = $d.Range($cell.Range.Start, $cell.Range.End)
$r.Find.Execute("0.6 (-3.8, 5.1)", $true, $false, $false, $false, $false, $true, 0, $false, "0.553 (-3.821, 5.125)", 2)
$cell = $t.Cell(3, 8)
$r = $d.Range($cell.Range.Start, $cell.Range.End)
$r.Find.Execute("12.4 (4.6, 20.7)*", $true, $false, $false, $false, $false, $true, 0, $false, "12.377 (4.624, 20.705)*", 2)
$cell = $t.Cell(3, 10)
$r = $d.Range($cell.Range.Start, $cell.Range.End)
$r.Find.Execute("30.9 (27.0, 34.8)*", $true, $false, $false, $false, $false, $true, 0, $false, "30.876 (27.046, 34.821)*", 2)
$cell = $t.Cell(4, 2)
$r = $d.Range($cell.Range.Start, $cell.Range.End)
$r.Find.Execute("0.8 (0.7, 0.9)", $true, $false, $false, $false, $false, $true, 0, $false, "0.833 (0.733, 0.933)", 2)
$cell = $t.Cell(4, 3)
$r = $d.Range($cell.Range.Start, $cell.Range.End)
$r.Find.Execute("3.1 (2.9, 3.3)", $true, $false, $false, $false, $false, $true, 0, $false, "3.111 (2.942, 3.279)", 2)
$cell = $t.Cell(4, 4)
$r = $d.Range($cell.Range.Start, $cell.Range.End)
$r.Find.Execute("9.0 (6.8, 11.2)*", $true, $false, $false, $false, $false, $true, 0, $false, "9.017 (6.832, 11.247)*", 2)
$cell = $t.Cell(4, 6)
$r = $d.Range($cell.Range.Start, $cell.Range.End)
$r.Find.Execute("-0.6 (-2.9, 1.8)", $true, $false, $false, $false, $false, $true, 0, $false, "-0.562 (-2.860, 1.790)", 2)
$cell = $t.Cell(4, 8)
$r = $d.Range($cell.Range.Start, $cell.Range.End)
$r.Find.Execute("33.5 (26.9, 40.4)*", $true, $false, $false, $false, $false, $true, 0, $false, "33.464 (26.903, 40.365)*", 2)
$cell = $t.Cell(6, 2)
$r = $d.Range($cell.Range.Start, $cell.Range.End)
$r.Find.Execute("0.3 (0.3, 0.3)", $true, $false, $false, $false, $false, $true, 0, $false, "0.313 (0.290, 0.336)", 2)
$cell = $t.Cell(6, 3)
$r = $d.Range($cell.Range.Start, $cell.Range.End)
$r.Find.Execute("1.2 (1.2, 1.2)", $true, $false, $false, $false, $false, $true, 0, $false, "1.204 (1.159, 1.248)", 2)
$cell = $t.Cell(6, 4)
$r = $d.Range($cell.Range.Start, $cell.Range.End)
$r.Find.Execute("9.2 (6.6, 11.8)*", $true, $false, $false, $false, $false, $true, 0, $false, "9.193 (6.625, 11.822)*", 2)
$cell = $t.Cell(6, 6)
$r = $d.Range($cell.Range.Start, $cell.Range.End)
$r.Find.Execute("42.7 (33.8, 52.3)*", $true, $false, $false, $false, $false, $true, 0, $false, "42.719 (33.779, 52.256)*", 2)
$cell = $t.Cell(6, 8)
$r = $d.Range($cell.Range.Start, $cell.Range.End)
$r.Find.Execute("21.4 (8.3, 36.1)*", $true, $false, $false, $false, $false, $true, 0, $false, "21.385 (8.297, 36.054)*", 2)
$cell = $t.Cell(6, 10)
$r = $d.Range($cell.Range.Start, $cell.Range.End)
$r.Find.Execute("-6.4 (-7.4, -5.4)*", $true, $false, $false, $false, $false, $true, 0, $false, "-6.419 (-7.400, -5.427)*", 2)
$cell = $t.Cell(7, 2)
$r = $d.Range($cell.Range.Start, $cell.Range.End)
$r.Find.Execute("0.2 (0.2, 0.3)", $true, $false, $false, $false, $false, $true, 0, $false, "0.213 (0.163, 0.263)", 2)
$cell = $t.Cell(7, 3)
$r = $d.Range($cell.Range.Start, $cell.Range.End)
$r.Find.Execute("0.6 (0.5, 0.7)", $true, $false, $false, $false, $false, $true, 0, $false, "0.578 (0.506, 0.651)", 2)
$cell = $t.Cell(7, 4)
$r = $d.Range($cell.Range.Start, $cell.Range.End)
$r.Find.Execute("6.2 (4.1, 8.3)*", $true, $false, $false, $false, $false, $true, 0, $false, "6.215 (4.138, 8.334)*", 2)
$cell = $t.Cell(7, 6)
$r = $d.Range($cell.Range.Start, $cell.Range.End)
$r.Find.Execute("21.8 (16.7, 27.1)*", $true, $false, $false, $false, $false, $true, 0, $false, "21.806 (16.692, 27.144)*", 2)
$cell = $t.Cell(7, 8)
$r = $d.Range($cell.Range.Start, $cell.Range.End)
$r.Find.Execute("-4.5 (-6.4, -2.6)*", $true, $false, $false, $false, $false, $true, 0, $false, "-4.518 (-6.434, -2.563)*", 2)
$cell = $t.Cell(9, 2)
$r = $d.Range($cell.Range.Start, $cell.Range.End)
$r.Find.Execute("1.1 (1.1, 1.2)", $true, $false, $false, $false, $false, $true, 0, $false, "1.109 (1.066, 1.153)", 2)
$cell = $t.Cell(9, 3)
$r = $d.Range($cell.Range.Start, $cell.Range.End)
$r.Find.Execute("4.6 (4.5, 4.7)", $true, $false, $false, $false, $false, $true, 0, $false, "4.601 (4.515, 4.687)", 2)
$cell = $t.Cell(9, 4)
$r = $d.Range($cell.Range.Start, $cell.Range.End)
$r.Find.Execute("8.9 (7.8, 10.0)*", $true, $false, $false, $false, $false, $true, 0, $false, "8.898 (7.799, 10.008)*", 2)
$cell = $t.Cell(9, 6)
$r = $d.Range($cell.Range.Start, $cell.Range.End)
$r.Find.Execute("12.8 (11.5, 14.1)*", $true, $false, $false, $false, $false, $true, 0, $false, "12.810 (11.511, 14.125)*", 2)
$cell = $t.Cell(9, 8)
$r = $d.Range($cell.Range.Start, $cell.Range.End)
$r.Find.Execute("0.8 (-1.8, 3.4)", $true, $false, $false, $false, $false, $true, 0, $false, "0.761 (-1.779, 3.367)", 2)
$cell = $t.Cell(10, 2)
$r = $d.Range($cell.Range.Start, $cell.Range.End)
$r.Find.Execute("0.6 (0.5, 0.7)", $true, $false, $false, $false, $false, $true, 0, $false, "0.622 (0.536, 0.709)", 2)
$cell = $t.Cell(10, 3)
$r = $d.Range($cell.Range.Start, $cell.Range.End)
$r.Find.Execute("2.1 (2.0, 2.3)", $true, $false, $false, $false, $false, $true, 0, $false, "2.110 (1.969, 2.250)", 2)
$cell = $t.Cell(10, 4)
$r = $d.Range($cell.Range.Start, $cell.Range.End)
$r.Find.Execute("8.2 (7.7, 8.7)*", $true, $false, $false, $false, $false, $true, 0, $false, "8.163 (7.658, 8.671)*", 2)
$cell = $t.Cell(10, 6)
$r = $d.Range($cell.Range.Start, $cell.Range.End)
$r.Find.Execute("8.2 (7.7, 8.7)*", $true, $false, $false, $false, $false, $true, 0, $false, "8.163 (7.658, 8.671)*", 2)
$cell = $t.Cell(12, 2)
$r = $d.Range($cell.Range.Start, $cell.Range.End)
$r.Find.Execute("0.3 (0.3, 0.3)", $true, $false, $false, $false, $false, $true, 0, $false, "0.298 (0.275, 0.320)", 2)
$cell = $t.Cell(12, 3)
$r = $d.Range($cell.Range.Start, $cell.Range.End)
$r.Find.Execute("3.6 (3.5, 3.7)", $true, $false, $false, $false, $false, $true, 0, $false, "3.585 (3.507, 3.662)", 2)
$cell = $t.Cell(12, 4)
$r = $d.Range($cell.Range.Start, $cell.Range.End)
$r.Find.Execute("17.1 (13.2, 21.1)*", $true, $false, $false, $false, $false, $true, 0, $false, "17.063 (13.196, 21.063)*", 2)
$cell = $t.Cell(12, 6)
$r = $d.Range($cell.Range.Start, $cell.Range.End)
$r.Find.Execute("17.7 (11.3, 24.4)*", $true, $false, $false, $false, $false, $true, 0, $false, "17.685 (11.296, 24.441)*", 2)
$cell = $t.Cell(12, 8)
$r = $d.Range($cell.Range.Start, $cell.Range.End)
$r.Find.Execute("3.2 (-1.7, 8.3)", $true, $false, $false, $false, $false, $true, 0, $false, "3.189 (-1.661, 8.278)", 2)
$cell = $t.Cell(12, 10)
$r = $d.Range($cell.Range.Start, $cell.Range.End)
$r.Find.Execute("78.7 (50.4, 112.3)*", $true, $false, $false, $false, $false, $true, 0, $false, "78.700 (50.424, 112.290)*", 2)
$cell = $t.Cell(13, 2)
$r = $d.Range($cell.Range.Start, $cell.Range.End)
$r.Find.Execute("0.1 (0.1, 0.2)", $true, $false, $false, $false, $false, $true, 0, $false, "0.115 (0.078, 0.152)", 2)
$cell = $t.Cell(13, 3)
$r = $d.Range($cell.Range.Start, $cell.Range.End)
$r.Find.Execute("2.1 (2.0, 2.2)", $true, $false, $false, $false, $false, $true, 0, $false, "2.098 (1.960, 2.237)", 2)
$cell = $t.Cell(13, 4)
$r = $d.Range($cell.Range.Start, $cell.Range.End)
$r.Find.Execute("15.4 (0.7, 32.2)*", $true, $false, $false, $false, $false, $true, 0, $false, "15.355 (0.674, 32.176)*", 2)
$cell = $t.Cell(13, 6)
$r = $d.Range($cell.Range.Start, $cell.Range.End)
$r.Find.Execute("6.1 (-5.7, 19.3)", $true, $false, $false, $false, $false, $true, 0, $false, "6.086 (-5.683, 19.323)", 2)
$cell = $t.Cell(13, 8)
$r = $d.Range($cell.Range.Start, $cell.Range.End)
$r.Find.Execute("107.3 (-14.7, 403.7)", $true, $false, $false, $false, $false, $true, 0, $false, "107.343 (-14.652, 403.719)", 2)
$cell = $t.Cell(15, 2)
$r = $d.Range($cell.Range.Start, $cell.Range.End)
$r.Find.Execute("1.0 (0.9, 1.0)", $true, $false, $false, $false, $false, $true, 0, $false, "0.964 (0.923, 1.004)", 2)
$cell = $t.Cell(15, 3)
$r = $d.Range($cell.Range.Start, $cell.Range.End)
$r.Find.Execute("1.0 (1.0, 1.0)", $true, $false, $false, $false, $false, $true, 0, $false, "1.005 (0.964, 1.046)", 2)
$cell = $t.Cell(15, 4)
$r = $d.Range($cell.Range.Start, $cell.Range.End)
$r.Find.Execute("-0.2 (-0.8, 0.4)", $true, $false, $false, $false, $false, $true, 0, $false, "-0.227 (-0.808, 0.357)", 2)
$cell = $t.Cell(15, 6)
$r = $d.Range($cell.Range.Start, $cell.Range.End)
$r.Find.Execute("-0.2 (-0.8, 0.4)", $true, $false, $false, $false, $false, $true, 0, $false, "-0.227 (-0.808, 0.357)", 2)
$cell = $t.Cell(16, 2)
$r = $d.Range($cell.Range.Start, $cell.Range.End)
$r.Find.Execute("2.1 (1.9, 2.2)", $true, $false, $false, $false, $false, $true, 0, $false, "2.056 (1.900, 2.212)", 2)
$cell = $t.Cell(16, 3)
$r = $d.Range($cell.Range.Start, $cell.Range.End)
$r.Find.Execute("0.6 (0.5, 0.7)", $true, $false, $false, $false, $false, $true, 0, $false, "0.593 (0.519, 0.667)", 2)
$cell = $t.Cell(16, 4)
$r = $d.Range($cell.Range.Start, $cell.Range.End)
$r.Find.Execute("-8.2 (-9.1, -7.3)*", $true, $false, $false, $false, $false, $true, 0, $false, "-8.215 (-9.133, -7.288)*", 2)
$cell = $t.Cell(16, 6)
$r = $d.Range($cell.Range.Start, $cell.Range.End)
$r.Find.Execute("-8.2 (-9.1, -7.3)*", $true, $false, $false, $false, $false, $true, 0, $false, "-8.215 (-9.133, -7.288)*", 2)

Write-Host "Done. Total replacements: 50"
